$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 733.3333
$ws.Range("J7").Value = 1500
$ws.Range("L7").Value = 1500
$ws.Range("N7").Value = -1724

$ws.Range("H14").Value = 733.3333
$ws.Range("J14").Value = 1500
$ws.Range("L14").Value = 1500
$ws.Range("N14").Value = -1882

$ws.Range("H21").Value = 44833.25

$ws.Range("H23").Value = 44833.25

$ws.Range("H32").Value = 9193.75
$ws.Range("I32").Value = 3760
$ws.Range("J32").Value = 18250
$ws.Range("K32").Value = 3760
$ws.Range("L32").Value = 18250
$ws.Range("M32").Value = -3434
$ws.Range("N32").Value = -18902

$ws.Range("H39").Value = 6251.3335
$ws.Range("I39").Value = 109
$ws.Range("J39").Value = 27749.5
$ws.Range("K39").Value = 327
$ws.Range("L39").Value = 83248.5
$ws.Range("M39").Value = -31
$ws.Range("N39").Value = -83840.5

$ws.Range("H46").Value = 3208.2144
$ws.Range("I46").Value = 975
$ws.Range("K46").Value = 2925
$ws.Range("M46").Value = -2806

$ws.Range("H60").Value = 3208.2144
$ws.Range("I60").Value = 975
$ws.Range("K60").Value = 2925
$ws.Range("M60").Value = -2441

$ws.Range("H87").Value = 368999.4
$ws.Range("J87").Value = 368999.4
$ws.Range("L87").Value = 368999.4
$ws.Range("N87").Value = -371495.4

$ws.Range("H90").Value = 368999.4
$ws.Range("J90").Value = 368999.4
$ws.Range("L90").Value = 1106998.2
$ws.Range("N90").Value = -1119478.2

$ws.Range("H101").Value = 1300.4615
$ws.Range("J101").Value = 696.25
$ws.Range("L101").Value = 2088.75
$ws.Range("N101").Value = -5332.75

$ws.Range("H132").Value = 7706.696
$ws.Range("I132").Value = 5695.1953
$ws.Range("K132").Value = 17085.5859
$ws.Range("M132").Value = -14555.5859

$ws.Range("H135").Value = 25754.572
$ws.Range("I135").Value = 8898.833000000001
$ws.Range("J135").Value = 48228.89
$ws.Range("K135").Value = 80089.497
$ws.Range("L135").Value = 434060.01
$ws.Range("M135").Value = -77554.497
$ws.Range("N135").Value = -439130.01

$ws.Range("H137").Value = 13399.171
$ws.Range("I137").Value = 5087.696
$ws.Range("K137").Value = 15263.088
$ws.Range("M137").Value = -12713.088

$ws.Range("H138").Value = 5489.755
$ws.Range("J138").Value = 5904.617
$ws.Range("L138").Value = 17713.851
$ws.Range("N138").Value = -27993.851

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 4200.25
$ws.Range("J97").Value = 6514.3
$ws.Range("L97").Value = 6514.3
$ws.Range("N97").Value = -7506.3

$ws.Range("H102").Value = 17979.928
$ws.Range("I102").Value = 3476.6667
$ws.Range("K102").Value = 3476.6667
$ws.Range("M102").Value = -1854.6667

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H137").Value = 59222.223
$ws.Range("J137").Value = 59000
$ws.Range("L137").Value = 59000
$ws.Range("N137").Value = -69200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4014.4
$ws.Range("I94").Value = 2693.4707
$ws.Range("K94").Value = 2693.4707
$ws.Range("M94").Value = -2242.4707

$ws.Range("H105").Value = 9401.6
$ws.Range("I105").Value = 9401.6
$ws.Range("K105").Value = 9401.6
$ws.Range("M105").Value = -7654.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14741.292
$ws.Range("I31").Value = 5719.4443
$ws.Range("K31").Value = 5719.4443
$ws.Range("M31").Value = -5424.4443

$ws.Range("H34").Value = 14741.292
$ws.Range("I34").Value = 5719.4443
$ws.Range("K34").Value = 5719.4443
$ws.Range("M34").Value = -5517.4443

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H135").Value = 82119.766
$ws.Range("J135").Value = 82119.766
$ws.Range("L135").Value = 82119.766
$ws.Range("N135").Value = -92259.766

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1425.8667
$ws.Range("J92").Value = 1443.4445
$ws.Range("L92").Value = 4330.333500000001
$ws.Range("N92").Value = -6826.333500000001

$ws.Range("H113").Value = 2454.8572
$ws.Range("J113").Value = 2486.9
$ws.Range("L113").Value = 7460.700000000001
$ws.Range("N113").Value = -11800.7

$ws.Range("H132").Value = 1757.92
$ws.Range("I132").Value = 1718.5714
$ws.Range("J132").Value = 1808
$ws.Range("K132").Value = 15467.1426
$ws.Range("L132").Value = 16272
$ws.Range("M132").Value = -12937.1426
$ws.Range("N132").Value = -21332

$ws.Range("H140").Value = 1056.1875
$ws.Range("I140").Value = 1056.1875
$ws.Range("K140").Value = 3168.5625
$ws.Range("M140").Value = 2011.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 4079.8
$ws.Range("I2").Value = 4079.8
$ws.Range("K2").Value = 4079.8
$ws.Range("M2").Value = -3966.8

$ws.Range("H3").Value = 2501425
$ws.Range("I3").Value = 3333566.8
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 3333566.8
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -3333450.8
$ws.Range("N3").Value = -5232

$ws.Range("H9").Value = 806.3333
$ws.Range("I9").Value = 806.3333
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 806.3333
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -636.3333
$ws.Range("N9").ClearContents()

$ws.Range("H40").Value = 70000
$ws.Range("I40").Value = 10000
$ws.Range("K40").Value = 10000
$ws.Range("M40").Value = -9849

$ws.Range("H46").Value = 46000
$ws.Range("J46").Value = 46000
$ws.Range("L46").Value = 46000
$ws.Range("N46").Value = -46312

$ws.Range("H80").Value = 18723
$ws.Range("I80").Value = 13562
$ws.Range("J80").Value = 23310.555
$ws.Range("K80").Value = 13562
$ws.Range("L80").Value = 23310.555
$ws.Range("M80").Value = -12564
$ws.Range("N80").Value = -25306.555

$ws.Range("H83").Value = 18723
$ws.Range("I83").Value = 13562
$ws.Range("J83").Value = 23310.555
$ws.Range("K83").Value = 67810
$ws.Range("L83").Value = 116552.775
$ws.Range("M83").Value = -62818
$ws.Range("N83").Value = -126536.775

$ws.Range("H122").Value = 4937.6665
$ws.Range("I122").Value = 2736.7896
$ws.Range("K122").Value = 8210.3688
$ws.Range("M122").Value = -5760.3688

$ws.Range("H133").Value = 82788.664
$ws.Range("J133").Value = 82788.664
$ws.Range("L133").Value = 82788.664
$ws.Range("N133").Value = -92908.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 2857.1428
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5224

$ws.Range("H15").Value = 2857.1428
$ws.Range("J15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("N15").Value = -5340

$ws.Range("H46").Value = 2007181.6
$ws.Range("J46").Value = 8964
$ws.Range("L46").Value = 8964
$ws.Range("N46").Value = -9340

$ws.Range("H68").Value = 6062.067
$ws.Range("J68").Value = 8249
$ws.Range("L68").Value = 8249
$ws.Range("N68").Value = -9747

$ws.Range("H69").Value = 80000
$ws.Range("J69").Value = 80000
$ws.Range("L69").Value = 80000
$ws.Range("N69").Value = -81622

$ws.Range("H71").Value = 6062.067
$ws.Range("J71").Value = 8249
$ws.Range("L71").Value = 41245
$ws.Range("N71").Value = -48733

$ws.Range("H72").Value = 80000
$ws.Range("J72").Value = 80000
$ws.Range("L72").Value = 240000
$ws.Range("N72").Value = -248112

$ws.Range("H93").Value = 8061.148
$ws.Range("I93").Value = 5150.421
$ws.Range("J93").Value = 14974.125
$ws.Range("K93").Value = 5150.421
$ws.Range("L93").Value = 14974.125
$ws.Range("M93").Value = -3902.421
$ws.Range("N93").Value = -17470.125

$ws.Range("H94").Value = 38900
$ws.Range("J94").Value = 38900
$ws.Range("L94").Value = 38900
$ws.Range("N94").Value = -40252

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 1600
$ws.Range("J54").Value = 2500
$ws.Range("L54").Value = 2500
$ws.Range("N54").Value = -3540

$ws.Range("H81").Value = 1831.5883
$ws.Range("I81").Value = 1831.5883
$ws.Range("K81").Value = 3663.1766
$ws.Range("M81").Value = -2602.1766

$ws.Range("H84").Value = 1831.5883
$ws.Range("I84").Value = 1831.5883
$ws.Range("K84").Value = 18315.883
$ws.Range("M84").Value = -13011.883
